# "add fun for save data and change method in bst"
#  - rename a station ("Новокузнецкая" -> "Новокосино") on the existing row
#  - append a new station row ("Новокузнецкая" with its own coordinates)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the station at row 2; its coordinates stay the same.
$ws.Range("A2").Value = "Новокосино"

# Append a brand new row with the same look & feel as the rows above it
# (copy row 5's formatting down into the new row 6, then overwrite the
# values with the new station's data).
$ws.Rows(5).Copy()
$ws.Rows(6).Insert()
$ws.Range("A6").Value = "Новокузнецкая"
$ws.Range("B6").Value = "55.74212,`n37.62901"
$ws.Rows(6).RowHeight = 24.4

[void]$ws.Range("A7").Select()
